$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename columns (dataset -> coded column names) ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Title-case the Spanish connector words (de/del/el/la/las/los/y) ---
# in state (col A) and municipality (col B) names throughout the data rows.
$nameFixes = @{
    'B8' = 'Pabellón De Arteaga'
    'B9' = 'Rincón De Romos'
    'B10' = 'San Francisco De Los Romo'
    'B11' = 'San José De Gracia'
    'B31' = 'Amatenango De La Frontera'
    'B33' = 'Benemérito De Las Américas'
    'B35' = 'Chiapa De Corzo'
    'B39' = 'Comitán De Domínguez'
    'B51' = 'Mazapa De Madero'
    'B59' = 'Salto De Agua'
    'B90' = 'Coyame Del Sotol'
    'B101' = 'Guadalupe Y Calvo'
    'B104' = 'Hidalgo Del Parral'
    'B128' = 'San Francisco De Borja'
    'B129' = 'San Francisco De Conchos'
    'B130' = 'San Francisco Del Oro'
    'B138' = 'Valle De Zaragoza'
    'B161' = 'San Juan De Sabinas'
    'B174' = 'Villa De Álvarez'
    'A176' = 'Ciudad De México'
    'B180' = 'Cuajimalpa De Morelos'
    'B194' = 'Coneto De Comonfort'
    'B208' = 'Nombre De Dios'
    'B212' = 'Pánuco De Coronado'
    'B219' = 'San Juan De Guadalupe'
    'B220' = 'San Juan Del Río'
    'B221' = 'San Luis Del Cordero'
    'B222' = 'San Pedro Del Gallo'
    'A232' = 'Estado De México'
    'B232' = 'Acambay De Ruíz Castañeda'
    'B235' = 'Almoloya De Juárez'
    'B240' = 'Atizapán De Zaragoza'
    'B251' = 'Ecatepec De Morelos'
    'B256' = 'Ixtapan De La Sal'
    'B265' = 'Naucalpan De Juárez'
    'B272' = 'San Felipe Del Progreso'
    'B274' = 'San Simón De Guerrero'
    'B283' = 'Tenango Del Valle'
    'B293' = 'Tlalnepantla De Baz'
    'B298' = 'Valle De Bravo'
    'B299' = 'Valle De Chalco Solidaridad'
    'B300' = 'Villa De Allende'
    'B301' = 'Villa Del Carbón'
    'B312' = 'San Miguel De Allende'
    'B313' = 'Apaseo El Alto'
    'B314' = 'Apaseo El Grande'
    'B322' = 'Dolores Hidalgo Cuna De La Independencia Nacional'
    'B326' = 'Jaral Del Progreso'
    'B334' = 'Purísima Del Rincón'
    'B338' = 'San Diego De La Unión'
    'B340' = 'San Francisco Del Rincón'
    'B342' = 'San Luis De La Paz'
    'B344' = 'Santa Cruz De Juventino Rosas'
    'B345' = 'Silao De La Victoria'
    'B349' = 'Valle De Santiago'
    'B354' = 'Acapulco De Juárez'
    'B356' = 'Ajuchitlán Del Progreso'
    'B357' = 'Alcozauca De Guerrero'
    'B360' = 'Atenango Del Río'
    'B362' = 'Atoyac De Álvarez'
    'B363' = 'Ayutla De Los Libres'
    'B366' = 'Buenavista De Cuéllar'
    'B367' = 'Chilapa De Álvarez'
    'B368' = 'Chilpancingo De Los Bravo'
    'B369' = 'Coahuayutla De José María Izazaga'
    'B374' = 'Coyuca De Benítez'
    'B375' = 'Coyuca De Catalán'
    'B379' = 'Cuetzala Del Progreso'
    'B380' = 'Cutzamala De Pinzón'
    'B386' = 'Huitzuco De Los Figueroa'
    'B387' = 'Iguala De La Independencia'
    'B389' = 'Ixcateopan De Cuauhtémoc'
    'B390' = 'Zihuatanejo De Azueta'
    'B392' = 'La Unión De Isidoro Montes De Oca'
    'B407' = 'Taxco De Alarcón'
    'B409' = 'Técpan De Galeana'
    'B411' = 'Tepecoacuilco De Trujano'
    'B413' = 'Tixtla De Guerrero'
    'B416' = 'Tlalixtaquilla De Maldonado'
    'B425' = 'Agua Blanca De Iturbide'
    'B429' = 'Atotonilco El Grande'
    'B434' = 'Cuautepec De Hinojosa'
    'B438' = 'Huejutla De Reyes'
    'B441' = 'Jacala De Ledezma'
    'B447' = 'Mineral Del Chico'
    'B448' = 'Mixquiahuala De Juárez'
    'B449' = 'Molango De Escamilla'
    'B450' = 'Nopala De Villagrán'
    'B451' = 'Pachuca De Soto'
    'B460' = 'Tepeji Del Río De Ocampo'
    'B461' = 'Tezontepec De Aldama'
    'B466' = 'Tula De Allende'
    'B467' = 'Tulancingo De Bravo'
    'B468' = 'Zacualtipán De Ángeles'
    'B472' = 'Acatlán De Juárez'
    'B473' = 'Ahualulco De Mercado'
    'B477' = 'Atemajac De Brizuela'
    'B480' = 'Atotonilco El Alto'
    'B482' = 'Autlán De Navarro'
    'B487' = 'Cañadas De Obregón'
    'B492' = 'Cuautitlán De García Barragán'
    'B499' = 'Encarnación De Díaz'
    'B503' = 'Huejuquilla El Alto'
    'B504' = 'Ixtlahuacán De Los Membrillos'
    'B505' = 'Ixtlahuacán Del Río'
    'B509' = 'Jilotlán De Los Dolores'
    'B513' = 'Lagos De Moreno'
    'B518' = 'Ojuelos De Jalisco'
    'B523' = 'San Diego De Alejandría'
    'B525' = 'San Juan De Los Lagos'
    'B526' = 'San Juanito De Escobedo'
    'B528' = 'San Martín De Bolaños'
    'B530' = 'San Miguel El Alto'
    'B531' = 'San Sebastián Del Oeste'
    'B532' = 'Santa María De Los Ángeles'
    'B535' = 'Talpa De Allende'
    'B536' = 'Tamazula De Gordiano'
    'B542' = 'Tepatitlán De Morelos'
    'B545' = 'Tizapán El Alto'
    'B546' = 'Tlajomulco De Zúñiga'
    'B554' = 'Unión De San Antonio'
    'B555' = 'Unión De Tula'
    'B556' = 'Valle De Guadalupe'
    'B561' = 'Yahualica De González Gallo'
    'B562' = 'Zacoalco De Torres'
    'B565' = 'Zapotlán Del Rey'
    'B566' = 'Zapotlán El Grande'
    'B589' = 'Cojumatlán De Régules'
    'B644' = 'Tiquicheo De Nicolás Romero'
    'B675' = 'Puente De Ixtla'
    'B679' = 'Tetela Del Volcán'
    'B680' = 'Tlaltizapán De Zapata'
    'B689' = 'Amatlán De Cañas'
    'B690' = 'Bahía De Banderas'
    'B693' = 'Ixtlán Del Río'
    'B700' = 'Santa María Del Oro'
    'B711' = 'Ciénega De Flores'
    'B717' = 'Mier Y Noriega'
    'B719' = 'San Nicolás De Los Garza'
    'B723' = 'Acatlán De Pérez Figueroa'
    'B725' = 'Ayoquezco De Aldama'
    'B727' = 'Chalcatongo De Hidalgo'
    'B729' = 'Fresnillo De Trujano'
    'B730' = 'Heroica Ciudad De Ejutla De Crespo'
    'B731' = 'Heroica Ciudad De Huajuapan De León'
    'B732' = 'Heroica Ciudad De Tlaxiaco'
    'B733' = 'Huautla De Jiménez'
    'B734' = 'Ixtlán De Juárez'
    'B735' = 'Heroica Ciudad De Juchitán De Zaragoza'
    'B738' = 'Mariscala De Juárez'
    'B741' = 'Miahuatlán De Porfirio Díaz'
    'B742' = 'Nejapa De Madero'
    'B743' = 'Oaxaca De Juárez'
    'B744' = 'Ocotlán De Morelos'
    'B745' = 'Putla Villa De Guerrero'
    'B764' = 'San José Del Peñasco'
    'B765' = 'San Juan Bautista Lo De Soto'
    'B768' = 'San Juan Del Estado'
    'B790' = 'San Pablo Villa De Mitla'
    'B798' = 'San Pedro Y San Pablo Ayutla'
    'B839' = 'Tataltepec De Valdés'
    'B840' = 'Heroica Villa Tezoatlán De Segura Y Luna, Cuna De La Independencia De Oaxaca'
    'B841' = 'Tlacolula De Matamoros'
    'B842' = 'Villa De Chilapa De Díaz'
    'B843' = 'Villa De Etla'
    'B844' = 'Villa De Tututepec'
    'B846' = 'Villa Sola De Vega'
    'B847' = 'Zimatlán De Álvarez'
    'B865' = 'Ayotoxco De Guerrero'
    'B869' = 'Chalchicomula De Sesma'
    'B884' = 'Huehuetlán El Chico'
    'B885' = 'Huehuetlán El Grande'
    'B888' = 'Huitzilan De Serdán'
    'B891' = 'Izúcar De Matamoros'
    'B900' = 'Los Reyes De Juárez'
    'B901' = 'Mazapiltepec De Juárez'
    'B906' = 'Palmar De Bravo'
    'B915' = 'San Nicolás De Los Ranchos'
    'B917' = 'San Salvador El Seco'
    'B922' = 'Tecali De Herrera'
    'B927' = 'Tepanco De López'
    'B928' = 'Tepatlaxco De Hidalgo'
    'B932' = 'Tepexi De Rodríguez'
    'B934' = 'Tetela De Ocampo'
    'B935' = 'Teteles De Avila Castillo'
    'B940' = 'Tlacotepec De Benito Juárez'
    'B958' = 'Amealco De Bonfil'
    'B960' = 'Cadereyta De Montes'
    'B966' = 'Jalpan De Serra'
    'B967' = 'Landa De Matamoros'
    'B970' = 'Pinal De Amoles'
    'B973' = 'San Juan Del Río'
    'B983' = 'Armadillo De Los Infante'
    'B984' = 'Axtla De Terrazas'
    'B995' = 'Mexquitic De Carmona'
    'B1001' = 'Santa María Del Río'
    'B1003' = 'Soledad De Graciano Sánchez'
    'B1008' = 'Tanquián De Escobedo'
    'B1011' = 'Villa De Arriaga'
    'B1012' = 'Villa De Ramos'
    'B1013' = 'Villa De Reyes'
    'B1054' = 'Nacozari De García'
    'B1092' = 'Soto La Marina'
    'B1103' = 'Contla De Juan Cuamatzi'
    'B1106' = 'Ixtacuixtla De Mariano Matamoros'
    'B1107' = 'Nanacamilpa De Mariano Arista'
    'B1110' = 'San Pablo Del Monte'
    'B1111' = 'Sanctórum De Lázaro Cárdenas'
    'B1114' = 'Tepetitla De Lardizábal'
    'B1117' = 'Tetla De La Solidaridad'
    'B1128' = 'Alto Lucero De Gutiérrez Barrios'
    'B1131' = 'Amatlán De Los Reyes'
    'B1147' = 'Cosamaloapan De Carpio'
    'B1158' = 'Hueyapan De Ocampo'
    'B1159' = 'Ignacio De La Llave'
    'B1161' = 'Ixhuatlán Del Café'
    'B1170' = 'Lerdo De Tejada'
    'B1173' = 'Martínez De La Torre'
    'B1186' = 'Paso De Ovejas'
    'B1187' = 'Paso Del Macho'
    'B1191' = 'Poza Rica De Hidalgo'
    'B1195' = 'Sayula De Alemán'
    'B1198' = 'Soledad De Doblado'
    'B1220' = 'Vega De Alatorre'
    'B1253' = 'Cañitas De Felipe Pescador'
    'B1255' = 'Concepción Del Oro'
    'B1257' = 'El Plateado De Joaquín Amaro'
    'B1266' = 'Jiménez Del Teul'
    'B1275' = 'Moyahua De Estrada'
    'B1276' = 'Nochistlán De Mejía'
    'B1277' = 'Noria De Ángeles'
    'B1288' = 'Teúl De González Ortega'
    'B1289' = 'Tlaltenango De Sánchez Román'
    'B1291' = 'Villa De Cos'
}
foreach ($cellRef in $nameFixes.Keys) {
    $ws.Range($cellRef).Value = $nameFixes[$cellRef]
}

# --- Recomputed percentage values (tiny float recalculation deltas) ---
$pctFixes = @{
    'D65' = 0.0009807767752059633
    'D74' = 0.009415457041977248
    'D76' = 0.0009807767752059633
    'D92' = 0.0009807767752059633
    'D94' = 0.0009153916568588988
    'D130' = 0.0009153916568588988
    'D132' = 0.0009807767752059633
    'D146' = 0.0009807767752059633
    'D188' = 0.0009153916568588988
    'D194' = 0.0009807767752059633
    'D206' = 0.0009807767752059633
    'D215' = 0.0009153916568588988
    'D272' = 0.0009153916568588988
    'D332' = 0.0009807767752059633
    'D340' = 0.0009807767752059633
    'D350' = 0.0009807767752059633
    'D358' = 0.0009807767752059633
    'D451' = 0.0009153916568588988
    'D510' = 0.0009153916568588988
    'D543' = 0.0009807767752059633
    'D547' = 0.0009153916568588988
    'D624' = 0.0009807767752059633
    'D659' = 0.0009153916568588988
    'D766' = 0.0009153916568588988
    'D945' = 0.0009807767752059633
    'D981' = 0.0009153916568588988
    'D988' = 0.0009153916568588988
    'D991' = 0.0009807767752059633
    'D1002' = 0.0009807767752059633
    'D1029' = 0.0009807767752059633
    'D1035' = 0.0009153916568588988
    'D1056' = 0.0009807767752059633
    'D1174' = 0.0009153916568588988
    'D1190' = 0.0009807767752059633
}
foreach ($cellRef in $pctFixes.Keys) {
    $ws.Range($cellRef).Value = $pctFixes[$cellRef]
}

# --- Drop the trailing footnote rows (source/sample-size notes) ---
$ws.Range("A1299:A1304").EntireRow.Delete()
